$wb = $excel.ActiveWorkbook

# Map of sheet index (1-based, matches xl/worksheets/sheetN.xml order)
# to the F-column ("想去人数") row updates described in the diff.

$ws = $wb.Worksheets.Item(1)
$updates = @{
    2 = 279
    3 = 956
    4 = 1225
    5 = 1074
    6 = 3173
    9 = 1134
    10 = 679
    11 = 48
    12 = 550
    13 = 272
    14 = 26
    16 = 1318
    17 = 1318
    18 = 223
    20 = 575
    21 = 261
    22 = 542
    23 = 428
    24 = 685
    26 = 309
    27 = 31
    28 = 882
    29 = 160
    31 = 386
    32 = 1103
    33 = 5180
    34 = 616
    35 = 352
    36 = 253
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item(2)
$updates = @{
    15 = 54
    16 = 316
    19 = 430
    22 = 754
    40 = 778
    42 = 53
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item(3)
$updates = @{
    4 = 676
    5 = 480
    6 = 472
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item(4)
$updates = @{
    2 = 676
    3 = 279
    4 = 480
    6 = 956
    7 = 1225
    9 = 1074
    10 = 3173
    14 = 1134
    15 = 679
    16 = 472
    18 = 48
    19 = 550
    20 = 272
    22 = 26
    24 = 1318
    25 = 1318
    26 = 224
    27 = 54
    28 = 575
    29 = 261
    30 = 428
    33 = 685
    35 = 309
    36 = 31
    37 = 882
    38 = 160
    40 = 386
    41 = 1103
    42 = 5180
    43 = 616
    45 = 352
    46 = 253
    48 = 778
    50 = 53
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()